$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Text corrections (label rewording)
# ---------------------------------------------------------------------------
$ws.Cells.Replace("Kali Pertama", "Semakan Kali Pertama")
$ws.Cells.Replace("Kali Kedua", "Semakan Kali Kedua")
$ws.Cells.Replace("Kali Ketiga", "Semakan Kali Ketiga")
$ws.Cells.Replace("Kali Keempat", "Semakan Kali Keempat")
$ws.Cells.Replace("BOUQUET KREATIF", "Bouquet Kreatif")
$ws.Cells.Replace("TIK TOK RAYA", "Tik Tok Raya")
$ws.Cells.Replace("RIANG RIA KUIH RAYA", "Riang Ria Kuih Raya")
$ws.Cells.Replace("CREATIVE COLLAGE", "Creative Collage")

# ---------------------------------------------------------------------------
# 2) Move the report title from E4 to D4 (merged heading now starts at D4)
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = $ws.Range("E4").Value2
$ws.Range("E4").ClearContents()
